# feat: add heat wave duration variable and update related translations and styles
#
# Adds a new "hw_dur" worksheet (heat-wave duration climatology) right after
# the existing "hw_int" (heat-wave intensity) sheet, fills it with the same
# header/label layout as hw_int (reusing the shared strings), applies a new
# "0.00_ " number format to the numeric data cells, and makes the new sheet
# the active tab/selection.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- create the new sheet right after hw_int -------------------------------
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$newSheet.Name = "hw_dur"

# --- header rows (same labels/layout as hw_int) -----------------------------
# Row 1: scenario names
$newSheet.Range("B1").Value = "Historical"
$newSheet.Range("C1").Value = "SSP2-4.5"
$newSheet.Range("E1").Value = "SSP3-7.0"
$newSheet.Range("G1").Value = "SSP5-8.5"

# Row 2: period labels
$newSheet.Range("C2").Value = "2046-2065"
$newSheet.Range("D2").Value = "2081-2100"
$newSheet.Range("E2").Value = "2046-2065"
$newSheet.Range("F2").Value = "2081-2100"
$newSheet.Range("G2").Value = "2046-2065"
$newSheet.Range("H2").Value = "2081-2100"

# Touch the blank header cells (A1, D1, F1, H1, A2, B2) so they are emitted
# as real (empty) cells, matching the formatted A1:H14 block used for the
# rest of the sheet.
$newSheet.Range("A1:H2").WrapText = $false

# --- monthly duration data (degC*days or similar, style = new 0.00_ format) -
$newSheet.Range("B3:H14").NumberFormat = "0.00_ "

$rowVals = @(3.11, 3.34, 3.45, 3.46, 3.64, 3.36, 4.19)
$newSheet.Range("A3").Value = "Jan"
for ($i = 0; $i -lt $rowVals.Length; $i++) { $newSheet.Cells.Item(3, 2 + $i).Value = $rowVals[$i] }
$rowVals = @(3.11, 3.27, 3.43, 3.44, 3.61, 3.33, 4.14)
$newSheet.Range("A4").Value = "Feb"
for ($i = 0; $i -lt $rowVals.Length; $i++) { $newSheet.Cells.Item(4, 2 + $i).Value = $rowVals[$i] }
$rowVals = @(3.17, 3.27, 3.45, 3.41, 3.46, 3.29, 3.81)
$newSheet.Range("A5").Value = "Mar"
for ($i = 0; $i -lt $rowVals.Length; $i++) { $newSheet.Cells.Item(5, 2 + $i).Value = $rowVals[$i] }
$rowVals = @(3.13, 3.23, 3.38, 3.37, 3.44, 3.27, 3.76)
$newSheet.Range("A6").Value = "Apr"
for ($i = 0; $i -lt $rowVals.Length; $i++) { $newSheet.Cells.Item(6, 2 + $i).Value = $rowVals[$i] }
$rowVals = @(3.12, 3.22, 3.36, 3.34, 3.71, 3.32, 3.76)
$newSheet.Range("A7").Value = "May"
for ($i = 0; $i -lt $rowVals.Length; $i++) { $newSheet.Cells.Item(7, 2 + $i).Value = $rowVals[$i] }
$rowVals = @(3.13, 3.43, 3.43, 3.49, 3.73, 3.49, 4.07)
$newSheet.Range("A8").Value = "Jun"
for ($i = 0; $i -lt $rowVals.Length; $i++) { $newSheet.Cells.Item(8, 2 + $i).Value = $rowVals[$i] }
$rowVals = @(3.12, 3.4, 3.38, 3.43, 3.65, 3.46, 4)
$newSheet.Range("A9").Value = "Jul"
for ($i = 0; $i -lt $rowVals.Length; $i++) { $newSheet.Cells.Item(9, 2 + $i).Value = $rowVals[$i] }
$rowVals = @(3.17, 3.49, 3.49, 3.46, 3.8, 4.07, 4.14)
$newSheet.Range("A10").Value = "Aug"
for ($i = 0; $i -lt $rowVals.Length; $i++) { $newSheet.Cells.Item(10, 2 + $i).Value = $rowVals[$i] }
$rowVals = @(3.13, 3.3, 3.34, 3.45, 3.52, 3.23, 3.84)
$newSheet.Range("A11").Value = "Sep"
for ($i = 0; $i -lt $rowVals.Length; $i++) { $newSheet.Cells.Item(11, 2 + $i).Value = $rowVals[$i] }
$rowVals = @(3.13, 3.27, 3.38, 3.46, 3.51, 3.26, 3.85)
$newSheet.Range("A12").Value = "Oct"
for ($i = 0; $i -lt $rowVals.Length; $i++) { $newSheet.Cells.Item(12, 2 + $i).Value = $rowVals[$i] }
$rowVals = @(3.16, 3.4, 3.39, 3.52, 3.56, 3.31, 3.95)
$newSheet.Range("A13").Value = "Nov"
for ($i = 0; $i -lt $rowVals.Length; $i++) { $newSheet.Cells.Item(13, 2 + $i).Value = $rowVals[$i] }
$rowVals = @(3.12, 3.33, 3.46, 3.44, 3.62, 3.34, 4.28)
$newSheet.Range("A14").Value = "Dec"
for ($i = 0; $i -lt $rowVals.Length; $i++) { $newSheet.Cells.Item(14, 2 + $i).Value = $rowVals[$i] }

# --- make the new sheet the active tab/selection ----------------------------
$newSheet.Activate()
$newSheet.Range("A1:H14").Select() | Out-Null

Write-Output "hw_dur sheet added"
